# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets
# to reflect newly scraped counts from the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 496
$ws1.Range("F4").Value = 1262
$ws1.Range("F6").Value = 14090
$ws1.Range("F7").Value = 15608
$ws1.Range("F11").Value = 188
$ws1.Range("F23").Value = 6100
$ws1.Range("F24").Value = 960
$ws1.Range("F25").Value = 1091
$ws1.Range("F26").Value = 5547
$ws1.Range("F28").Value = 140
$ws1.Range("F29").Value = 114
$ws1.Range("F30").Value = 4517

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 496
$ws4.Range("F4").Value = 1262
$ws4.Range("F6").Value = 14090
$ws4.Range("F7").Value = 15608
$ws4.Range("F11").Value = 188
$ws4.Range("F24").Value = 6100
$ws4.Range("F25").Value = 960
$ws4.Range("F26").Value = 1091
$ws4.Range("F27").Value = 5547
$ws4.Range("F29").Value = 140
$ws4.Range("F30").Value = 114
$ws4.Range("F31").Value = 4517
